# Adam23-Itga5 LR-pair sheet: refresh values with the updated NATMI run
# (per Dr Hou's advice) -- adds "M2" as a 4th cluster and rewrites rows 2-17
# with the recomputed per-cluster ligand/receptor statistics and edge weights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: SendingCluster, Ligand, Receptor, TargetCluster,
#             E..J (ligand-side stats), K..P (receptor-side stats), Q..T (edge weights)
$rows = @(
    ("ECs", "Adam23", "Itga5", "ECs", 3, 1, 0.1448543333333333, 0.434563, 0.007514552731137001, 0.007514552731137, 3, 1, 31.22896466666667, 93.686894, 0.2877106972998646, 0.2877106972998646, 4.523650857480222, 40.712857717322, 0.002162017206172028, 0.002162017206172028),
    ("ECs", "Adam23", "Itga5", "FAPs", 3, 1, 0.1448543333333333, 0.434563, 0.007514552731137001, 0.007514552731137, 3, 1, 40.44578266666667, 121.337348, 0.3726247238124506, 0.3726247238124505, 5.858746884324889, 52.728721958924, 0.002800108136014021, 0.00280010813601402),
    ("ECs", "Adam23", "Itga5", "M2", 3, 1, 0.1448543333333333, 0.434563, 0.007514552731137001, 0.007514552731137, 3, 1, 25.36964133333333, 76.108924, 0.2337290805561598, 0.2337290805561598, 3.674902482245778, 33.074122340212, 0.00175636950063943, 0.00175636950063943),
    ("ECs", "Adam23", "Itga5", "sCs", 3, 1, 0.1448543333333333, 0.434563, 0.007514552731137001, 0.007514552731137, 3, 1, 11.49855033333333, 34.495651, 0.1059354983315251, 0.1059354983315251, 1.665614842834778, 14.990533585513, 0.0007960578883115214, 0.0007960578883115213),
    ("FAPs", "Adam23", "Itga5", "ECs", 3, 1, 14.494489, 43.483467, 0.7519250504625467, 0.7519250504625465, 3, 1, 31.22896466666667, 93.686894, 0.2877106972998646, 0.2877106972998646, 452.6478848423887, 4073.830963581498, 0.2163368805858152, 0.2163368805858151),
    ("FAPs", "Adam23", "Itga5", "FAPs", 3, 1, 14.494489, 43.483467, 0.7519250504625467, 0.7519250504625465, 3, 1, 40.44578266666667, 121.337348, 0.3726247238124506, 0.3726247238124505, 586.2409519583907, 5276.168567625517, 0.2801858642562694, 0.2801858642562693),
    ("FAPs", "Adam23", "Itga5", "M2", 3, 1, 14.494489, 43.483467, 0.7519250504625467, 0.7519250504625465, 3, 1, 25.36964133333333, 76.108924, 0.2337290805561598, 0.2337290805561598, 367.7199872399454, 3309.479885159508, 0.175746750691755, 0.175746750691755),
    ("FAPs", "Adam23", "Itga5", "sCs", 3, 1, 14.494489, 43.483467, 0.7519250504625467, 0.7519250504625465, 3, 1, 11.49855033333333, 34.495651, 0.1059354983315251, 0.1059354983315251, 166.6656113224464, 1499.990501902017, 0.07965555492870707, 0.07965555492870706),
    ("M2", "Adam23", "Itga5", "ECs", 2, 0.6666666666666666, 0.091182, 0.273546, 0.004730213666123443, 0.004730213666123443, 3, 1, 31.22896466666667, 93.686894, 0.2877106972998646, 0.2877106972998646, 2.847519456236, 25.627675106124, 0.001360933072257725, 0.001360933072257725),
    ("M2", "Adam23", "Itga5", "FAPs", 2, 0.6666666666666666, 0.091182, 0.273546, 0.004730213666123443, 0.004730213666123443, 3, 1, 40.44578266666667, 121.337348, 0.3726247238124506, 0.3726247238124505, 3.687927355112, 33.191346196008, 0.001762594560913128, 0.001762594560913127),
    ("M2", "Adam23", "Itga5", "M2", 2, 0.6666666666666666, 0.091182, 0.273546, 0.004730213666123443, 0.004730213666123443, 3, 1, 25.36964133333333, 76.108924, 0.2337290805561598, 0.2337290805561598, 2.313254636056, 20.819291724504, 0.001105588491017214, 0.001105588491017214),
    ("M2", "Adam23", "Itga5", "sCs", 2, 0.6666666666666666, 0.091182, 0.273546, 0.004730213666123443, 0.004730213666123443, 3, 1, 11.49855033333333, 34.495651, 0.1059354983315251, 0.1059354983315251, 1.048460816494, 9.436147348446001, 0.0005010975419353775, 0.0005010975419353775),
    ("sCs", "Adam23", "Itga5", "ECs", 3, 1, 4.545982333333334, 13.637947, 0.235830183140193, 0.2358301831401929, 3, 1, 31.22896466666667, 93.686894, 0.2877106972998646, 0.2877106972998646, 141.9663216629576, 1277.696894966618, 0.06785086643561969, 0.06785086643561968),
    ("sCs", "Adam23", "Itga5", "FAPs", 3, 1, 4.545982333333334, 13.637947, 0.235830183140193, 0.2358301831401929, 3, 1, 40.44578266666667, 121.337348, 0.3726247238124506, 0.3726247238124505, 183.8658134605062, 1654.792321144556, 0.08787615685925405, 0.08787615685925401),
    ("sCs", "Adam23", "Itga5", "M2", 3, 1, 4.545982333333334, 13.637947, 0.235830183140193, 0.2358301831401929, 3, 1, 25.36964133333333, 76.108924, 0.2337290805561598, 0.2337290805561598, 115.3299413043365, 1037.969471739028, 0.05512037187274807, 0.05512037187274806),
    ("sCs", "Adam23", "Itga5", "sCs", 3, 1, 4.545982333333334, 13.637947, 0.235830183140193, 0.2358301831401929, 3, 1, 11.49855033333333, 34.495651, 0.1059354983315251, 0.1059354983315251, 52.27220667427746, 470.4498600684971, 0.02498278797257118, 0.02498278797257118)
)

# Clear the old data rows (A2:T13) before laying the new (larger) table down.
$ws.Range("A2:T13").ClearContents()

$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value2 = $row[$i]
    }
    $r++
}
